$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Factor change (J3: 90 -> 100); dependent formulas in J4:K9 recalc automatically
$ws.Range("J3").Value2 = 100

# --- Fill in row 8 with a new "Implementacion" time entry (2021-10-19)
$ws.Range("B8").Value2 = "Implementacion"
$ws.Range("C8").Value2 = "2021-10-19"
$ws.Range("D8").Value2 = 0.60416666666666663
$ws.Range("E8").Value2 = 0.72916666666666663
$ws.Range("F8").Formula = "=E8-D8"
$ws.Range("G8").Formula = "=(HOUR(F8)*60)+MINUTE(F8)"
$ws.Range("H8").Formula = "=(H`$3/60)*(I`$3/100)"
$ws.Range("I8").Formula = "=G8*H8"
$ws.Range("J8").Formula = "=I8*J`$3%"
$ws.Range("K8").Formula = "=I8-J8"

# --- Fill in row 9 with another new "Implementacion" time entry (2021-10-19)
$ws.Range("B9").Value2 = "Implementacion"
$ws.Range("C9").Value2 = "2021-10-19"
$ws.Range("D9").Value2 = 0.8125
$ws.Range("E9").Value2 = 0.88680555555555562
$ws.Range("F9").Formula = "=E9-D9"
$ws.Range("G9").Formula = "=(HOUR(F9)*60)+MINUTE(F9)"
$ws.Range("H9").Formula = "=(H`$3/60)*(I`$3/100)"
$ws.Range("I9").Formula = "=G9*H9"
$ws.Range("J9").Formula = "=I9*J`$3%"
$ws.Range("K9").Formula = "=I9-J9"

# --- Remove the now-unneeded blank template rows (old rows 12-23)
$ws.Rows("12:23").Delete()

# --- Move the selection (cursor) to match the saved view
$ws.Range("B11").Select()
